$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the policy number (NumPoliza) in E2 - keep TIPO_ENDOSO (F2) unchanged
# Prefix with an apostrophe so Excel stores it as text (preserving the leading
# zero) without altering the cell's existing quote-prefix text style/format.
$ws.Range("E2").Value = "'04104016708"

# Move the active selection to E2 (was F3)
$ws.Range("E2").Select()
